# Word COM-interop script applying the documented edits.
#
# Approach:
#  - Plain text-only edits (no run split needed) are done with
#    $d.Content.Find.Execute(...).
#  - Edits that split an existing run into multiple runs (matching the
#    target OOXML) are done by locating the exact character offsets with
#    Find, then forcing a run boundary by toggling a character formatting
#    property (Font.Bold on/off, net no-op) on the sub-range - Word (and
#    this host) always materializes distinct <w:r> elements at range
#    boundaries when a direct formatting write touches only part of a run.
#  - The new list item is added with Range.InsertParagraphAfter(), which
#    inherits the paragraph/list formatting of the preceding item.

$d = $word.ActiveDocument

function Split-RunAt($rangeStart, $rangeEnd) {
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "Every thing that goes in the Database HAS TO BE AN OBJECT ..."
#    Split "Every thing" off from the rest of the first run.
# ---------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("Every thing that goes in the Database HAS TO BE AN OBJECT")
$pStart = $full.Start
$pMid = $pStart + 11
Split-RunAt $pStart $pMid

# Split "etter(" off the front of the final run ("etter() should exist...")
$full = $d.Content
$full.Find.Execute("etter() should exist for the value inserted.")
$pStart = $full.Start
$pMid = $pStart + 6
Split-RunAt $pStart $pMid

# ---------------------------------------------------------------------
# 2) "Every getter should return STRING or Boolean only." ->
#    "Every getter should return STRING or Boolean or Custom OBJECT only."
#    split across three runs.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Every getter should return STRING or Boolean only.", $true, $false, $false, $false, $false, $true, 1, $false, "Every getter should return STRING or Boolean or Custom OBJECT only.", 2)

$full = $d.Content
$full.Find.Execute("Every getter should return STRING or Boolean or Custom OBJECT only.")
$pStart = $full.Start
$pMid1 = $pStart + 45
$pMid2 = $pMid1 + 17
Split-RunAt $pStart $pMid1
Split-RunAt $pMid1 $pMid2

# ---------------------------------------------------------------------
# 3) "Your .gitignore needs to be ..." split into "Your .", "gitignore",
#    " needs to be".
# ---------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("Your .gitignore needs to be")
$pStart = $full.Start
$pMid1 = $pStart + 6
$pMid2 = $pMid1 + 9
Split-RunAt $pStart $pMid1
Split-RunAt $pMid1 $pMid2

# ---------------------------------------------------------------------
# 4) New list item after "Only the final version of the PR will be
#    merged." - "FORMAT – [Issue #xx] Branch name/something else."
# ---------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("Only the final version of the PR will be merged.")
$full.Collapse(0)
$full.InsertParagraphAfter()

$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Only the final version of the PR will be merged.*") {
        $targetIdx = $idx + 1
    }
}
$newPara = $d.Paragraphs($targetIdx)
$newPara.Range.Text = "FORMAT – [Issue #xx] Branch name/something else."
